$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column B ("Skill Description") before the existing SFIA Level column,
# shifting SFIA Level, Keycode, Description right by one column.
$ws.Columns.Item(2).Insert()

$ws.Range("B1").Value = "Skill Description"

$descriptions = @{
    2  = "Autonomy"
    3  = "Autonomy"
    4  = "Autonomy"
    5  = "Influence"
    6  = "Influence"
    7  = "Complexity"
    8  = "Complexity"
    9  = "Complexity"
    10 = "Knowledge"
    11 = "Knowledge"
    12 = "Acceptance testing"
    13 = "Acceptance testing"
    14 = "Testing"
    15 = "Testing"
    16 = "Testing"
    17 = "Testing"
    18 = "User experience evaluation"
    19 = "User experience evaluation"
    20 = "MADE"
}

foreach ($row in $descriptions.Keys) {
    $ws.Cells.Item($row, 2).Value = $descriptions[$row]
}
